$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add new label cells (literal text; reuses existing shared strings where
#     the text already exists elsewhere in the sheet, and creates a brand new
#     shared string for the new "annual_F_biomass_target" label) -------------

# Row 25: add C25 = "annual_F_biomass_target" (new string)
$ws.Range("C25").Value = "annual_F_biomass_target"

# Row 38: add A38 = "DERIVED_QUANTITIES" (module name, same as A2..A8)
$ws.Range("A38").Value = "DERIVED_QUANTITIES"

# Row 56/57: add A56/A57 = "OVERALL_COMPS" (module name, same as A9)
$ws.Range("A56").Value = "OVERALL_COMPS"
$ws.Range("A57").Value = "OVERALL_COMPS"

# Row 65,69,72,74: add A column = "CATCH" (module name, same as A66/A70/A73)
$ws.Range("A65").Value = "CATCH"
$ws.Range("A69").Value = "CATCH"
$ws.Range("A72").Value = "CATCH"
$ws.Range("A74").Value = "CATCH"

# Row 75: add C75 = "catch_weight" (same alt_label as row 70)
$ws.Range("C75").Value = "catch_weight"

# --- Column widths -----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 22.81640625
$ws.Columns.Item(2).ColumnWidth = 19.08984375
$ws.Columns.Item(3).ColumnWidth = 34.90625

# --- Sheet view: scroll position + selection --------------------------------
$excel.ActiveWindow.ScrollRow = 97
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C76").Select()
